# Remove the obsolete "Straftaten" (Z16_B01_P01 / indicator 16.1) row.
# Deleting the entire row 64 shifts every following row up by one,
# renumbers the sheet's `r` attributes automatically, and Excel keeps
# the worksheet's used-range dimension in sync (A1:O69 -> A1:O68).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(64).Delete()
